$d = $word.ActiveDocument

# 1. Collapse the "ngày @Day tháng @Month năm @Year" placeholder into a single "@Date" placeholder,
#    keeping the leading ", " text in its own run (matching the target OOXML).
$d.Content.Find.Execute("@Day tháng @Month năm @Year", $false, $false, $false, $false, $false, $true, 1, $false, "@Date", 2)
$d.Content.Find.Execute(", ngày ", $false, $false, $false, $false, $false, $true, 1, $false, ", ", 2)

# Nudge formatting on the new "@Date" run and back so it does not get silently
# re-coalesced into the preceding ", " run when the document is serialized.
$r = $d.Content.Duplicate
$r.Find.Execute("@Date", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# 2. Collapse "@Co_organizedUnit" into "@CoUnit".
$d.Content.Find.Execute("@Co_organizedUnit", $false, $false, $false, $false, $false, $true, 1, $false, "@CoUnit", 2)

# 3. Collapse "@AttendanceStart - @AttendanceEnd" into "@Attendance".
$d.Content.Find.Execute("@AttendanceStart - @AttendanceEnd", $false, $false, $false, $false, $false, $true, 1, $false, "@Attendance", 2)
